# China Southern RPA done
# Replace the WO/Reference/Waybill data rows with the new ChinaSouthernCargo data set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (headers) is unchanged: Waybill Number | WONumber | ReferenceNumber

# New data rows (A: Waybill Number, B: WONumber, C: ReferenceNumber)
$data = @(
    @("784-69784750", 2042807795,    2042807795),
    @("784-24772392", "T040384489",  "T040384489"),
    @("784-40756752", 2052898015,    2052898015),
    @("784-22861764", "DJMFEA4228904", "23G0010440"),
    @("784-69735175", "DJORDA4228886", 2042806833),
    @("784-69735175", "DJORDA4228898", 2042806834),
    @("784-69735061", "T210010236",  "T210010236"),
    @("784-22862416", "DJAMSA4227503", 2482456744),
    @("784-69785866", "DJAMSA4225276", 2482458634),
    @("784-69734884", "T260053631",  "T260053631")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
}

$wb.Save()
